# "Deep dive total run."
# Adds a second block (rows 5-10) to Sheet2 with a wider "deep dive" breakdown
# of the award/idv counts, mirroring the existing header/data block in A1:K3
# but with more columns (A:R) and a Comma[0]-style numeric format on the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ---- header row (row 5) ----
$ws.Range("A5").Value = "award_or_idv_flag"
$ws.Range("B5").Value = "IsDerived"
$ws.Range("C5").Value = "HasParent"
$ws.Range("D5").Value = "cauofferblank"
$ws.Range("E5").Value = "pcauofferblank"
$ws.Range("F5").Value = "ppcauofferblank"
$ws.Range("G5").Value = "cauidvblank"
$ws.Range("H5").Value = "pcauidvblank"
$ws.Range("I5").Value = "ppcauidvblank"
$ws.Range("J5").Value = "cauoidcblank"
$ws.Range("K5").Value = "pcauidcblank"
$ws.Range("L5").Value = "ppcauidcblank"
$ws.Range("M5").Value = "cauomultiblank"
$ws.Range("N5").Value = "pcaumultiblank"
$ws.Range("O5").Value = "ppcaumultiblank"
$ws.Range("P5").Value = "caucount"
$ws.Range("Q5").Value = "pcaucount"
$ws.Range("R5").Value = "nonderivedpcaucount"

# ---- data rows (6-10) ----
$data = @(
    @("AWARD", 0, 0, 4058,    4058,    4058,    2322539,  2322539, 2322539, 2322539,  2322539, 2322539, 2322539, 2322539, 2322539, 2322539,  0,        0),
    @("AWARD", 0, 1, 8264182, 4598638, 4478569, 12966376, 7479678, 7346447, 12966376, 9109364, 9109364, 226314,  211370,  211048,  12966376, 12966376, 7734502),
    @("IDV",   0, 0, 55612,   55612,   55612,   38334,    38334,   38334,   74164,    74164,   74164,   38422,   38422,   38422,   130813,   0,        0),
    @("IDV",   0, 1, 2326,    559,     559,     2307,     528,     528,     7034,     7034,    7034,    2310,    528,     528,     7034,     6991,     6521),
    @("IDV",   1, 0, 60917,   60917,   60917,   60917,    60917,   60917,   60917,    60917,   60917,   60917,   60917,   60917,   60917,    0,        0)
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = 6 + $i
    $rowVals = $data[$i]
    $ws.Range("A$r").Value = $rowVals[0]
    for ($j = 1; $j -lt $cols.Count; $j++) {
        $cell = $ws.Range($cols[$j] + "$r")
        $cell.Value = $rowVals[$j]
        $cell.Style = "Comma"
        $cell.NumberFormat = "_(* #,##0_);_(* \(#,##0\);_(* ""-""??_);_(@_)"
    }
}

# ---- selection mirrors what Excel leaves after the paste/fill ----
$ws.Range("A5:R10").Select()
